$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo "IINYO" -> "INYO" in the CB point-name strings (code review feedback).
$ws.Range("B4").Value = "INYO 115KV CB"
$ws.Range("B7").Value = "HAIWEE-INYOKERN 115KV CB"
$ws.Range("B8").Value = "COSO-HAIWEE-INYOKERN 115KV CB"
